$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Furniture section: insert a new "Drawer" row right after "Closet" (row 60),
#    pushing "Sofa"/"Lamp" (and everything below, including the "Job" block)
#    down by one row.
# ---------------------------------------------------------------------------
$ws.Rows("61").Insert()
$ws.Range("B61").Value = "Drawer"
$ws.Range("C61").Value = 0.34

# ---------------------------------------------------------------------------
# 2. New "Clothing" category block, typed in at rows 72-78 (row 71 left blank
#    as a separator, matching the sheet's existing layout convention), then
#    sorted by proportion descending - exactly like the other categories.
# ---------------------------------------------------------------------------
$ws.Range("A72").Value = "Clothing"

$items = @("Dress", "Shirt", "Pants", "Sock", "Shoe", "Jacket", "Jeans")
$props = @(0.56000000000000005, 0.27, 0.67, 0.7, 0.82, 0.93, 0.61)
for ($i = 0; $i -lt $items.Length; $i++) {
    $ws.Cells.Item(72 + $i, 2).Value = $items[$i]
    $ws.Cells.Item(72 + $i, 3).Value = $props[$i]
}

$sortRange = $ws.Range("B72:C78")
$sortKey = $ws.Range("C72:C78")
$sortRange.Sort($sortKey, 2)

# Match the look of the other category-header rows: a handful of decorative
# blank cells (re-using the existing "filler" style) and the taller row
# height that goes with them.
$ws.Range("D57").Copy()
$ws.Range("D72").PasteSpecial(-4122)
$ws.Range("G72").PasteSpecial(-4122)
$ws.Range("H72").PasteSpecial(-4122)
$ws.Range("I72").PasteSpecial(-4122)
$ws.Range("A72").EntireRow.RowHeight = 18

# ---------------------------------------------------------------------------
# 3. Highlight the category-header cells in column A with a yellow fill.
#    First application defines the new fill/style; later ones reuse it via a
#    format-only copy/paste so no duplicate styles get created.
# ---------------------------------------------------------------------------
$ws.Range("A3").Interior.Color = 65535
$ws.Range("A3").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A72").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Restore the view: scroll/select so that D67 is the active cell (matches
#    where the author was last looking while reorganizing the sheet).
# ---------------------------------------------------------------------------
$ws.Range("D67").Select()
